$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 8313
$ws.Range("F5").Value = 8313
$ws.Range("F8").Value = 117
$ws.Range("F9").Value = 7281
$ws.Range("F10").Value = 1145
$ws.Range("F11").Value = 575
$ws.Range("F19").Value = 126
$ws.Range("F20").Value = 12003
$ws.Range("F23").Value = 2411
$ws.Range("F24").Value = 3416
$ws.Range("F27").Value = 2852
$ws.Range("F29").Value = 33
$ws.Range("F30").Value = 3324
$ws.Range("F33").Value = 1683
$ws.Range("F34").Value = 78
$ws.Range("I34").Value = "//i0.hdslb.com/bfs/openplatform/202407/pKdcyAR31721272661076.jpeg"
$ws.Range("F35").Value = 119
$ws.Range("F36").Value = 5948
$ws.Range("F37").Value = 94
$ws.Range("F38").Value = 1820
$ws.Range("F40").Value = 24
$ws.Range("F41").Value = 884
$ws.Range("F48").Value = 1569
$ws.Range("F49").Value = 13
$ws.Range("F50").Value = 112

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F12").Value = 3
$ws.Range("F15").Value = 7
$ws.Range("F20").Value = 71

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 302
$ws.Range("F3").Value = 443
$ws.Range("F4").Value = 10

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F5").Value = 302
$ws.Range("F8").Value = 8313
$ws.Range("F11").Value = 117
$ws.Range("F12").Value = 7281
$ws.Range("F13").Value = 7281
$ws.Range("F14").Value = 575
$ws.Range("F21").Value = 126
$ws.Range("F23").Value = 12003
$ws.Range("F27").Value = 2411
$ws.Range("F28").Value = 2411
$ws.Range("F29").Value = 3416
$ws.Range("F31").Value = 33
$ws.Range("F32").Value = 7
$ws.Range("F33").Value = 3324
$ws.Range("F36").Value = 1683
$ws.Range("F37").Value = 78
$ws.Range("I37").Value = "//i0.hdslb.com/bfs/openplatform/202407/pKdcyAR31721272661076.jpeg"
$ws.Range("F38").Value = 119
$ws.Range("F39").Value = 5948
$ws.Range("F40").Value = 71
$ws.Range("F41").Value = 94
$ws.Range("F42").Value = 1820
$ws.Range("F45").Value = 24
$ws.Range("F46").Value = 884
$ws.Range("F52").Value = 112
